$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 39/40: coin swap (TheGraph <-> PEPE) ---
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"

# --- Column D (Price) updates ---
$ws.Range("D2").Value = "67.107.87"
$ws.Range("D3").Value = "3.454.26"
$ws.Range("D5").Value = "'578.37"
$ws.Range("D6").Value = "'186.96"
$ws.Range("D8").Value = "3.447.86"
$ws.Range("D11").Value = "'0.643"
$ws.Range("D12").Value = "'57.66"
$ws.Range("D14").Value = "'9.47"
$ws.Range("D15").Value = "4.001.74"
$ws.Range("D17").Value = "3.449.16"
$ws.Range("D18").Value = "67.065.72"
$ws.Range("D20").Value = "'12.04"
$ws.Range("D22").Value = "'481.47"
$ws.Range("D23").Value = "'17.49"
$ws.Range("D24").Value = "'5.38"
$ws.Range("D25").Value = "'4.33"
$ws.Range("D26").Value = "'89.38"
$ws.Range("D28").Value = "'10.93"
$ws.Range("D29").Value = "'9.00"
$ws.Range("D30").Value = "'31.24"
$ws.Range("D31").Value = "'7.33"
$ws.Range("D32").Value = "'603.61"
$ws.Range("D33").Value = "'64.85"
$ws.Range("D34").Value = "'11.80"
$ws.Range("D35").Value = "'0.113"
$ws.Range("D38").Value = "'36.85"
$ws.Range("D39").Value = "0.0₃0777"
$ws.Range("D40").Value = "'0.385"
$ws.Range("D41").Value = "'3.46"
$ws.Range("D42").Value = "3.196.05"
$ws.Range("D45").Value = "'2.55"
$ws.Range("D46").Value = "'3.25"
$ws.Range("D48").Value = "'2.68"
$ws.Range("D50").Value = "'8.65"
$ws.Range("D51").Value = "'3.19"

# --- Column E (Volume 1h) updates ---
$ws.Range("E2").Value = "  +2.18%  "
$ws.Range("E3").Value = "  +1.62%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E5").Value = "  +3.17%  "
$ws.Range("E6").Value = "  +6.04%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("E8").Value = "  +1.76%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  -1.24%  "
$ws.Range("E11").Value = "  +0.64%  "
$ws.Range("E12").Value = "  +6.81%  "
$ws.Range("E13").Value = "  -1.64%  "
$ws.Range("E14").Value = "  +2.50%  "
$ws.Range("E15").Value = "  +1.74%  "
$ws.Range("E16").Value = "  +2.88%  "
$ws.Range("E17").Value = "  +1.54%  "
$ws.Range("E18").Value = "  +2.17%  "
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("E20").Value = "  +1.06%  "
$ws.Range("E21").Value = "  +1.75%  "
$ws.Range("E22").Value = "  +4.48%  "
$ws.Range("E23").Value = "  +23.61%  "
$ws.Range("E24").Value = "  +9.46%  "
$ws.Range("E25").Value = "  +4.52%  "
$ws.Range("E26").Value = "  +2.07%  "
$ws.Range("E27").Value = "  +0.95%  "
$ws.Range("E28").Value = "  +1.75%  "
$ws.Range("E29").Value = "  +2.58%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E31").Value = "  +11.15%  "
$ws.Range("E32").Value = "  +3.73%  "
$ws.Range("E33").Value = "  +1.99%  "
$ws.Range("E34").Value = "  +2.26%  "
$ws.Range("E35").Value = "  +3.79%  "
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("E37").Value = "  +2.18%  "
$ws.Range("E38").Value = "  +2.06%  "
$ws.Range("E39").Value = "  +3.88%  "
$ws.Range("E40").Value = "  +2.54%  "
$ws.Range("E41").Value = "  -4.09%  "
$ws.Range("E42").Value = "  +2.88%  "
$ws.Range("E43").Value = "  +2.85%  "
$ws.Range("E44").Value = "  +2.39%  "
$ws.Range("E45").Value = "  +3.88%  "
$ws.Range("E46").Value = "  +2.10%  "
$ws.Range("E47").Value = "  +0.83%  "
$ws.Range("E48").Value = "  +15.99%  "
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("E50").Value = "  +2.69%  "
$ws.Range("E51").Value = "  +2.43%  "
